# Auto-generated Excel COM-interop edit script.
# Applies updated currentAveragePrice / LevePrice / LeveProfit figures
# (columns H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW worksheets,
# matching the "chore: update Sheets via scheduled runner" data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 479.16666
$ws.Range("I61").Value = 301.66666
$ws.Range("J61").Value = 656.6667
$ws.Range("K61").Value = 904.9999799999999
$ws.Range("L61").Value = 1970.0001
$ws.Range("M61").Value = -732.9999799999999
$ws.Range("N61").Value = -2314.0001
$ws.Range("H62").Value = 3404
$ws.Range("I62").Value = 3404
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3404
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -2780
$ws.Range("H65").Value = 3404
$ws.Range("I65").Value = 3404
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 17020
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -13900
$ws.Range("H68").Value = 20000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 20000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 20000
$ws.Range("N68").Value = -21498
$ws.Range("H71").Value = 20000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 20000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 60000
$ws.Range("N71").Value = -67488
$ws.Range("H109").Value = 30000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 30000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774
$ws.Range("H115").Value = 1124
$ws.Range("I115").Value = 655
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 1965
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = -398
$ws.Range("N115").Value = -12134
$ws.Range("H141").Value = 2707.4482
$ws.Range("I141").Value = 1317.5555
$ws.Range("J141").Value = 4981.8184
$ws.Range("K141").Value = 3952.6665
$ws.Range("L141").Value = 14945.4552
$ws.Range("M141").Value = 1227.3335
$ws.Range("N141").Value = -25305.4552

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1187.4138
$ws.Range("I2").Value = 1050.9166
$ws.Range("J2").Value = 1842.6
$ws.Range("K2").Value = 1050.9166
$ws.Range("L2").Value = 1842.6
$ws.Range("M2").Value = -937.9166
$ws.Range("N2").Value = -2068.6
$ws.Range("H32").Value = 4307.67
$ws.Range("I32").Value = 3306.032
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 3306.032
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -3019.032
$ws.Range("N32").Value = -20574
$ws.Range("H45").Value = 1197810.4
$ws.Range("I45").Value = 1516612.4
$ws.Range("J45").Value = 2303
$ws.Range("K45").Value = 1516612.4
$ws.Range("L45").Value = 2303
$ws.Range("M45").Value = -1516235.4
$ws.Range("N45").Value = -3057
$ws.Range("H97").Value = 6021.3887
$ws.Range("I97").Value = 9375.362999999999
$ws.Range("J97").Value = 750.8570999999999
$ws.Range("K97").Value = 9375.362999999999
$ws.Range("L97").Value = 750.8570999999999
$ws.Range("M97").Value = -8879.362999999999
$ws.Range("N97").Value = -1742.8571
$ws.Range("H116").Value = 1187.4138
$ws.Range("I116").Value = 1050.9166
$ws.Range("J116").Value = 1842.6
$ws.Range("K116").Value = 1050.9166
$ws.Range("L116").Value = 1842.6
$ws.Range("M116").Value = 1243.0834
$ws.Range("N116").Value = -6430.6
$ws.Range("H132").Value = 3847637.2
$ws.Range("I132").Value = 5103383.5
$ws.Range("J132").Value = 1914.875
$ws.Range("K132").Value = 15310150.5
$ws.Range("L132").Value = 5744.625
$ws.Range("M132").Value = -15307620.5
$ws.Range("N132").Value = -10804.625

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1187.4138
$ws.Range("I3").Value = 1050.9166
$ws.Range("J3").Value = 1842.6
$ws.Range("K3").Value = 1050.9166
$ws.Range("L3").Value = 1842.6
$ws.Range("M3").Value = -936.9166
$ws.Range("N3").Value = -2070.6
$ws.Range("H86").Value = 17859496
$ws.Range("I86").Value = 2169.4348
$ws.Range("J86").Value = 100003200
$ws.Range("K86").Value = 2169.4348
$ws.Range("L86").Value = 100003200
$ws.Range("M86").Value = -1046.4348
$ws.Range("N86").Value = -100005446
$ws.Range("H89").Value = 17859496
$ws.Range("I89").Value = 2169.4348
$ws.Range("J89").Value = 100003200
$ws.Range("K89").Value = 10847.174
$ws.Range("L89").Value = 500016000
$ws.Range("M89").Value = -5231.173999999999
$ws.Range("N89").Value = -500027232
$ws.Range("H105").Value = 3961.8572
$ws.Range("I105").Value = 3020
$ws.Range("J105").Value = 4485.1113
$ws.Range("K105").Value = 3020
$ws.Range("L105").Value = 4485.1113
$ws.Range("M105").Value = -1273
$ws.Range("N105").Value = -7979.1113

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2313.484
$ws.Range("I62").Value = 2229
$ws.Range("J62").Value = 2467.0908
$ws.Range("K62").Value = 2229
$ws.Range("L62").Value = 2467.0908
$ws.Range("M62").Value = -1605
$ws.Range("N62").Value = -3715.0908
$ws.Range("H65").Value = 2313.484
$ws.Range("I65").Value = 2229
$ws.Range("J65").Value = 2467.0908
$ws.Range("K65").Value = 11145
$ws.Range("L65").Value = 12335.454
$ws.Range("M65").Value = -8025
$ws.Range("N65").Value = -18575.454

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 27.473684
$ws.Range("I12").Value = 27.88889
$ws.Range("J12").Value = 27.1
$ws.Range("K12").Value = 83.66667
$ws.Range("L12").Value = 81.30000000000001
$ws.Range("M12").Value = 89.33333
$ws.Range("N12").Value = -427.3
$ws.Range("H39").Value = 501.16278
$ws.Range("I39").Value = 500
$ws.Range("J39").Value = 501.1905
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 1503.5715
$ws.Range("M39").Value = -1206
$ws.Range("N39").Value = -2091.5715
$ws.Range("H113").Value = 43479044
$ws.Range("I113").Value = 125000580
$ws.Range("J113").Value = 892.6667
$ws.Range("K113").Value = 375001740
$ws.Range("L113").Value = 2678.0001
$ws.Range("M113").Value = -374999570
$ws.Range("N113").Value = -7018.0001
$ws.Range("H119").Value = 10597.5
$ws.Range("I119").Value = 4528.6665
$ws.Range("J119").Value = 16666.334
$ws.Range("K119").Value = 13585.9995
$ws.Range("L119").Value = 49999.00199999999
$ws.Range("M119").Value = -8747.999500000002
$ws.Range("N119").Value = -59675.00199999999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11869.125
$ws.Range("I70").Value = 25163.334
$ws.Range("J70").Value = 3892.6
$ws.Range("K70").Value = 25163.334
$ws.Range("L70").Value = 3892.6
$ws.Range("M70").Value = -24893.334
$ws.Range("N70").Value = -4432.6
$ws.Range("H73").Value = 11869.125
$ws.Range("I73").Value = 25163.334
$ws.Range("J73").Value = 3892.6
$ws.Range("K73").Value = 25163.334
$ws.Range("L73").Value = 3892.6
$ws.Range("M73").Value = -24227.334
$ws.Range("N73").Value = -5764.6
$ws.Range("H80").Value = 11113908
$ws.Range("I80").Value = 22224616
$ws.Range("J80").Value = 3199.9333
$ws.Range("K80").Value = 22224616
$ws.Range("L80").Value = 3199.9333
$ws.Range("M80").Value = -22223618
$ws.Range("N80").Value = -5195.933300000001
$ws.Range("H83").Value = 11113908
$ws.Range("I83").Value = 22224616
$ws.Range("J83").Value = 3199.9333
$ws.Range("K83").Value = 111123080
$ws.Range("L83").Value = 15999.6665
$ws.Range("M83").Value = -111118088
$ws.Range("N83").Value = -25983.6665
$ws.Range("H126").Value = 3720.8948
$ws.Range("I126").Value = 2137.75
$ws.Range("J126").Value = 4872.273
$ws.Range("K126").Value = 6413.25
$ws.Range("L126").Value = 14616.819
$ws.Range("M126").Value = -3943.25
$ws.Range("N126").Value = -19556.819

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5416.4194
$ws.Range("I7").Value = 5635.2666
$ws.Range("J7").Value = 5211.25
$ws.Range("K7").Value = 5635.2666
$ws.Range("L7").Value = 5211.25
$ws.Range("M7").Value = -5523.2666
$ws.Range("N7").Value = -5435.25
$ws.Range("H25").Value = 2700
$ws.Range("I25").Value = 2700
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 2700
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -2470
$ws.Range("H40").Value = 4546.8623
$ws.Range("I40").Value = 5734.3076
$ws.Range("J40").Value = 3582.0625
$ws.Range("K40").Value = 5734.3076
$ws.Range("L40").Value = 3582.0625
$ws.Range("M40").Value = -5598.3076
$ws.Range("N40").Value = -3854.0625
$ws.Range("H122").Value = 7044.6523
$ws.Range("I122").Value = 5775.8184
$ws.Range("J122").Value = 8207.75
$ws.Range("K122").Value = 17327.4552
$ws.Range("L122").Value = 24623.25
$ws.Range("M122").Value = -14877.4552
$ws.Range("N122").Value = -29523.25
$ws.Range("H126").Value = 5416.4194
$ws.Range("I126").Value = 5635.2666
$ws.Range("J126").Value = 5211.25
$ws.Range("K126").Value = 16905.7998
$ws.Range("L126").Value = 15633.75
$ws.Range("M126").Value = -14435.7998
$ws.Range("N126").Value = -20573.75

Write-Output "Applied Ultima_Profits refresh."
